$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Format the new "branch" column cells (C15:C18) as text, right-aligned,
# matching the existing style used by the rest of column C.
$branchRange = $ws.Range("C15:C18")
$branchRange.NumberFormat = "@"
$branchRange.HorizontalAlignment = -4152

# Append the new module users (Cheque Cancelatorio, Inventario Permanente,
# ONP, Fallas de Dispositivos) with their branch codes to the Users sheet.
$ws.Cells.Item(16, 3).Value = "026"
$ws.Cells.Item(15, 1).Value = "F00273"
$ws.Cells.Item(15, 3).Value = "073"
$ws.Cells.Item(16, 1).Value = "F00644"
$ws.Cells.Item(17, 1).Value = "F01106"
$ws.Cells.Item(17, 3).Value = "006"
$ws.Cells.Item(18, 1).Value = "F02547"
$ws.Cells.Item(18, 3).Value = "089"

$ws.Range("G7").Select()
